# Update "想去人数" (want-to-go count) figures for the 草莓动漫节 and
# 第一届ANE·DACG动漫嘉年华 entries, on both the "展览" sheet and the
# "全部类型" sheet, reflecting newly scraped totals.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1389
$wsExhibit.Range("F5").Value = 661

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1389
$wsAll.Range("F6").Value = 661
